$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44
$ws.Range('A44').Value = 111410460
$ws.Range('B44').Value = 56543
$ws.Range('D44').Value = 'NT'
$ws.Range('E44').Value = 103021
$ws.Range('F44').Value = 'Talltita'
$ws.Range('G44').Value = 'Poecile montanus'
$ws.Range('H44').Value = '(Conrad von Baldenstein, 1827)'
$ws.Range('J44').ClearContents() | Out-Null
$ws.Range('M44').Value = 'obs i häcktid, lämplig biotop'
$ws.Range('P44').Value = 'Bennarby, Dannemora, Upl 479 m NW, Upl'
$ws.Range('Q44').Value = 652960.6254034473
$ws.Range('R44').Value = 6675623.861893108
$ws.Range('S44').Value = 74
$ws.Range('Z44').Value = '15:40'
$ws.Range('AB44').Value = '15:40'

# Row 45
$ws.Range('A45').Value = 111410479
$ws.Range('I45').Value = '''4'
$ws.Range('P45').Value = 'Bennarby, Dannemora, Upl 533 m NW, Upl'
$ws.Range('Q45').Value = 652957.9861328325
$ws.Range('R45').Value = 6675686.990895226
$ws.Range('S45').Value = 4
$ws.Range('Z45').Value = '19:30'
$ws.Range('AB45').Value = '19:34'

# Row 47
$ws.Range('A47').Value = 111410478
$ws.Range('B47').Value = 96348
$ws.Range('D47').Value = 'VU'
$ws.Range('E47').Value = 220787
$ws.Range('F47').Value = 'Knärot'
$ws.Range('G47').Value = 'Goodyera repens'
$ws.Range('H47').Value = '(L.) R. Br.'
$ws.Range('I47').Value = '''2'
$ws.Range('J47').Value = 'plantor/tuvor'
$ws.Range('M47').ClearContents() | Out-Null
$ws.Range('P47').Value = 'Bennarby, Dannemora, Upl 261 m NW, Upl'
$ws.Range('Q47').Value = 653010.7491410983
$ws.Range('R47').Value = 6675367.51151045
$ws.Range('S47').Value = 9
$ws.Range('Z47').Value = '15:21'
$ws.Range('AB47').Value = '15:22'

# Row 52
$ws.Range('A52').Value = 111410470
$ws.Range('B52').Value = 95535
$ws.Range('D52').Value = 'LC'
$ws.Range('E52').Value = 221946
$ws.Range('F52').Value = 'Mattlummer'
$ws.Range('G52').Value = 'Lycopodium clavatum'
$ws.Range('H52').Value = 'L.'
$ws.Range('I52').ClearContents() | Out-Null  # -> empty text
$ws.Range('J52').ClearContents() | Out-Null  # -> empty text
$ws.Range('P52').Value = 'Bennarby, Dannemora, Upl 432 m NW, Upl'
$ws.Range('Q52').Value = 652955.775129037
$ws.Range('R52').Value = 6675560.91665418
$ws.Range('Z52').Value = '15:31'
$ws.Range('AB52').Value = '15:31'

# Row 53
$ws.Range('A53').Value = 111410475
$ws.Range('I53').Value = '''1'
$ws.Range('J53').Value = 'm²'
$ws.Range('P53').Value = 'Bennarby, Dannemora, Upl 590 m NW, Upl'
$ws.Range('Q53').Value = 653059.7969692796
$ws.Range('R53').Value = 6675792.827983578
$ws.Range('Z53').Value = '15:59'
$ws.Range('AB53').Value = '15:59'

# Row 54
$ws.Range('A54').Value = 111448644
$ws.Range('B54').Value = 103288
$ws.Range('E54').Value = 221144
$ws.Range('F54').Value = 'Grönpyrola'
$ws.Range('G54').Value = 'Pyrola chlorantha'
$ws.Range('H54').Value = 'Sw.'
$ws.Range('I54').Value = '''10'
$ws.Range('J54').Value = 'plantor/tuvor'
# NOTE: L54 should become an empty placeholder cell (engine cannot represent empty-string distinct from blank); leaving blank.
$ws.Range('P54').Value = 'Bennarby 405 m NW, Upl'
$ws.Range('Q54').Value = 652956.8549192698
$ws.Range('R54').Value = 6675547.018880062
$ws.Range('Z54').Value = '15:38'
$ws.Range('AB54').Value = '15:39'
$ws.Range('AF54').ClearContents() | Out-Null

# Row 55
$ws.Range('A55').Value = 111448641
$ws.Range('B55').Value = 89802
$ws.Range('E55').Value = 5420
$ws.Range('F55').Value = 'Grovticka'
$ws.Range('G55').Value = 'Phaeolus schweinitzii'
$ws.Range('H55').Value = '(Fr.) Pat.'
$ws.Range('I55').Value = '''1'
$ws.Range('J55').ClearContents() | Out-Null  # -> empty text
$ws.Range('L55').ClearContents() | Out-Null
$ws.Range('P55').Value = 'Bennarby 248 m NW, Upl'
$ws.Range('Q55').Value = 652992.5167218229
$ws.Range('R55').Value = 6675302.508242311
$ws.Range('Z55').Value = '13:30'
$ws.Range('AB55').Value = '13:30'
# NOTE: AF55 should become an empty placeholder cell (engine cannot represent empty-string distinct from blank); leaving blank.
